# "break out stock.yaml completed"
# - Convert D100:D107 (bsecode) on the "day" sheet from text to numeric.
# - Append 9 new stock rows (108-116) to the "day" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- Fix D100:D107: these were stored as text, should be plain numbers ---
$ws.Cells.Item(100, 4).Value = 500034
$ws.Cells.Item(101, 4).Value = 542652
$ws.Cells.Item(102, 4).Value = 500480
$ws.Cells.Item(103, 4).Value = 517354
$ws.Cells.Item(104, 4).Value = 500086
$ws.Cells.Item(105, 4).Value = 533155
$ws.Cells.Item(106, 4).Value = 500547
$ws.Cells.Item(107, 4).Value = 500477

# --- Append new rows 108-116 ---
# Columns: A=sr, B=nsecode, C=name, D=bsecode, E=per_chg, F=close, G=volume, H=timeframe, I=Date Time
$newRows = @(
    @(108, 1, "PAGEIND",    "Page Industries Limited",              "532827", 0.02,  39425,    11878,    "day", "10/07/2024 11:36:39"),
    @(109, 2, "COLPAL",     "Colgate Palmolive (india) Limited",    "500830", 1.5,   3040.5,   596437,   "day", "10/07/2024 11:36:39"),
    @(110, 3, "ASIANPAINT", "Asian Paints Limited",                  "500820", 3.15,  2996.45,  2299950,  "day", "10/07/2024 11:36:39"),
    @(111, 4, "SHRIRAMFIN", "Shriram Finance Ltd",                   "511218", -1.07, 2762.1,   1660147,  "day", "10/07/2024 11:36:39"),
    @(112, 5, "UBL",        "United Breweries Limited",              "532478", 0.24,  2110,     111020,   "day", "10/07/2024 11:36:39"),
    @(113, 6, "DALBHARAT",  "Dalmia Bharat Limited",                 "533309", -0.12, 1865.3,   477669,   "day", "10/07/2024 11:36:39"),
    @(114, 7, "AUBANK",     "AU Small Finance Bank",                 "540611", -1.58, 629.95,   6039883,  "day", "10/07/2024 11:36:39"),
    @(115, 8, "INDHOTEL",   "The Indian Hotels Company Limited",     "500850", 0.14,  609.85,   1927077,  "day", "10/07/2024 11:36:39"),
    @(116, 9, "VEDL",       "Vedanta Limited",                       "500295", -1.92, 456.7,    11556309, "day", "10/07/2024 11:36:39")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    # bsecode (column D) stays text-typed even though it looks numeric -
    # force text storage with a leading quote, same as typing it in Excel.
    $ws.Cells.Item($r, 4).Value = "'" + $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
}
